$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Periodo Mora" / "Valor Mora" table occupies rows 16-26 (E & F columns).
# This edit reverses the order of those rows: the previous EC (estado de
# cuenta) periods are removed and new ones are added in their place, which
# in effect mirrors the 11-row block top-to-bottom.

$startRow = 16
$endRow = 26

# Capture the original values before we start overwriting anything.
$periodVals = @()
$valorVals = @()
for ($r = $startRow; $r -le $endRow; $r++) {
    $periodVals += ,$ws.Range("E$r").Value()
    $valorVals  += ,$ws.Range("F$r").Value()
}

$rowCount = $endRow - $startRow + 1

for ($i = 0; $i -lt $rowCount; $i++) {
    $srcIndex = $rowCount - 1 - $i
    $destRow = $startRow + $i
    $ws.Range("E$destRow").Value = $periodVals[$srcIndex]
    $ws.Range("F$destRow").Value = $valorVals[$srcIndex]
}
